$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "69.818.67"
Set-TextValue $ws "E2" "  +0.05%  "

Set-TextValue $ws "D3" "3.788.93"
Set-TextValue $ws "E3" "  -0.66%  "

Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  +0.05%  "

Set-TextValue $ws "D5" "657.98"
Set-TextValue $ws "E5" "  +4.00%  "

Set-TextValue $ws "D6" "166.75"
Set-TextValue $ws "E6" "  +0.91%  "

Set-TextValue $ws "D7" "3.790.62"
Set-TextValue $ws "E7" "  -0.55%  "

Set-TextValue $ws "E8" "  -0.03%  "

Set-TextValue $ws "E9" "  +1.20%  "

Set-TextValue $ws "E10" "  -1.24%  "

Set-TextValue $ws "D11" "0.457"
Set-TextValue $ws "E11" "  +0.69%  "

Set-TextValue $ws "D12" "6.94"
Set-TextValue $ws "E12" "  +3.74%  "

Set-TextValue $ws "E13" "  -3.71%  "

Set-TextValue $ws "D14" "35.16"
Set-TextValue $ws "E14" "  -2.22%  "

Set-TextValue $ws "D15" "4.417.33"
Set-TextValue $ws "E15" "  -1.09%  "

Set-TextValue $ws "D16" "3.776.08"
Set-TextValue $ws "E16" "  +1.11%  "

Set-TextValue $ws "D17" "69.872.28"
Set-TextValue $ws "E17" "  +0.02%  "

Set-TextValue $ws "D18" "17.71"
Set-TextValue $ws "E18" "  -1.99%  "

Set-TextValue $ws "E19" "  +0.18%  "

Set-TextValue $ws "D20" "7.06"
Set-TextValue $ws "E20" "  -1.43%  "

Set-TextValue $ws "D21" "471.27"
Set-TextValue $ws "E21" "  +0.16%  "

Set-TextValue $ws "D22" "9.65"
Set-TextValue $ws "E22" "  -0.72%  "

Set-TextValue $ws "D23" "0.712"
Set-TextValue $ws "E23" "  +0.57%  "

Set-TextValue $ws "D24" "0.0000145"
Set-TextValue $ws "E24" "  -3.96%  "

Set-TextValue $ws "D25" "82.38"
Set-TextValue $ws "E25" "  -1.77%  "

Set-TextValue $ws "D26" "12.30"
Set-TextValue $ws "E26" "  +0.87%  "

Set-TextValue $ws "D27" "10.35"
Set-TextValue $ws "E27" "  +2.66%  "

Set-TextValue $ws "E28" "  -2.45%  "

Set-TextValue $ws "E29" "  +0.17%  "

Set-TextValue $ws "D30" "3.939.56"
Set-TextValue $ws "E30" "  -0.89%  "

Set-TextValue $ws "D31" "2.76"
Set-TextValue $ws "E31" "  +2.10%  "

Set-TextValue $ws "E32" "  +2.09%  "

Set-TextValue $ws "D33" "7.25"
Set-TextValue $ws "E33" "  -1.36%  "

Set-TextValue $ws "E34" "  -1.25%  "

Set-TextValue $ws "D35" "0.175"
Set-TextValue $ws "E35" "  +15.68%  "

Set-TextValue $ws "E36" "  -0.13%  "

Set-TextValue $ws "D37" "3.745.57"
Set-TextValue $ws "E37" "  -0.50%  "

Set-TextValue $ws "D38" "8.91"
Set-TextValue $ws "E38" "  -1.70%  "

Set-TextValue $ws "E39" "  -1.85%  "

Set-TextValue $ws "D40" "3.31"
Set-TextValue $ws "E40" "  -1.29%  "

Set-TextValue $ws "D41" "5.87"
Set-TextValue $ws "E41" "  -0.52%  "

Set-TextValue $ws "B42" "Mantle"
Set-TextValue $ws "C42" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D42" "0.963"
Set-TextValue $ws "E42" "  -1.63%  "

Set-TextValue $ws "B43" "FirstDigitalUSD"
Set-TextValue $ws "C43" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D43" "0.998"
Set-TextValue $ws "E43" "  -0.28%  "

Set-TextValue $ws "D45" "46.09"
Set-TextValue $ws "E45" "  +6.62%  "

Set-TextValue $ws "E46" "  +3.08%  "

Set-TextValue $ws "D47" "157.75"
Set-TextValue $ws "E47" "  +1.49%  "

Set-TextValue $ws "D48" "47.95"
Set-TextValue $ws "E48" "  +1.37%  "

Set-TextValue $ws "E49" "  -0.36%  "

Set-TextValue $ws "D50" "1.41"
Set-TextValue $ws "E50" "  +0.47%  "

Set-TextValue $ws "D51" "8.46"
Set-TextValue $ws "E51" "  -0.12%  "
